# Update column F ("dSF") values per repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 2
    3  = 1
    4  = 2
    5  = -1
    6  = -1
    7  = -3
    8  = 5
    9  = 3
    10 = -2
    11 = -1
    12 = -2
    13 = -3
    14 = -4
    15 = -1
    16 = 5
    17 = 4
    18 = 10
    19 = -4
    20 = 6
    21 = 1
    22 = 4
    23 = -1
    24 = 4
    26 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
